$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text, matching the
# original inline-string cell type (avoids Excel auto-converting to a number).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.353.65'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.56'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '286.26'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3748'
$ws.Range('E7').Value = '  +2.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3267'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.38'
$ws.Range('E9').Value = '  -5.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.139'
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07393'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.826'
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.794'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.573.09'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001094'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06717'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.79'
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.333'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.23'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.66'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.360.74'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.298'
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('E26').Value = '  -3.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.74'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.38'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.903'
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.96'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.747.23'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.050'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.898'
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.902'
$ws.Range('E34').Value = '  -4.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.490'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08205'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.283'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06285'
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2177'
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.243'
$ws.Range('E41').Value = '  -3.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.98'
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6075'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.74'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.739'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5889'
$ws.Range('E47').Value = '  -2.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.991'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.66'
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.177'
$ws.Range('E50').Value = '  -3.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07130'
$ws.Range('E51').Value = '  -1.20%  '
